$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 6210.7
$ws.Range("J86").Value = 5763
$ws.Range("L86").Value = 5763
$ws.Range("N86").Value = -8009
$ws.Range("H88").Value = 18938.555
$ws.Range("I88").Value = 44466.152
$ws.Range("J88").Value = 4509.913
$ws.Range("K88").Value = 44466.152
$ws.Range("L88").Value = 4509.913
$ws.Range("M88").Value = -44060.152
$ws.Range("N88").Value = -5321.913
$ws.Range("H89").Value = 6210.7
$ws.Range("J89").Value = 5763
$ws.Range("L89").Value = 28815
$ws.Range("N89").Value = -40047
$ws.Range("H91").Value = 18938.555
$ws.Range("I91").Value = 44466.152
$ws.Range("J91").Value = 4509.913
$ws.Range("K91").Value = 44466.152
$ws.Range("L91").Value = 4509.913
$ws.Range("M91").Value = -43062.152
$ws.Range("N91").Value = -7317.913
$ws.Range("H107").Value = 732.44446
$ws.Range("I107").Value = 706.63635
$ws.Range("J107").Value = 773
$ws.Range("K107").Value = 706.63635
$ws.Range("L107").Value = 773
$ws.Range("M107").Value = 1213.36365
$ws.Range("N107").Value = -4613
$ws.Range("H118").Value = 539.1667
$ws.Range("I118").Value = 369.0625
$ws.Range("K118").Value = 1107.1875
$ws.Range("M118").Value = 549.8125
$ws.Range("H132").Value = 1446.7858
$ws.Range("I132").Value = 1446.7858
$ws.Range("K132").Value = 4340.357400000001
$ws.Range("M132").Value = -1810.357400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7136.1025
$ws.Range("I32").Value = 4723.4688
$ws.Range("J32").Value = 18165.285
$ws.Range("K32").Value = 4723.4688
$ws.Range("L32").Value = 18165.285
$ws.Range("M32").Value = -4436.4688
$ws.Range("N32").Value = -18739.285
$ws.Range("H74").Value = 45638.652
$ws.Range("I74").Value = 56975.777
$ws.Range("J74").Value = 4825
$ws.Range("K74").Value = 56975.777
$ws.Range("L74").Value = 4825
$ws.Range("M74").Value = -56101.777
$ws.Range("N74").Value = -6573
$ws.Range("H77").Value = 45638.652
$ws.Range("I77").Value = 56975.777
$ws.Range("J77").Value = 4825
$ws.Range("K77").Value = 284878.885
$ws.Range("L77").Value = 24125
$ws.Range("M77").Value = -280510.885
$ws.Range("N77").Value = -32861
$ws.Range("H88").Value = 1722.7778
$ws.Range("J88").Value = 1824.75
$ws.Range("L88").Value = 1824.75
$ws.Range("N88").Value = -2636.75
$ws.Range("H91").Value = 1722.7778
$ws.Range("J91").Value = 1824.75
$ws.Range("L91").Value = 1824.75
$ws.Range("N91").Value = -4632.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 647.0278
$ws.Range("I80").Value = 959.17645
$ws.Range("J80").Value = 367.73685
$ws.Range("K80").Value = 959.17645
$ws.Range("L80").Value = 367.73685
$ws.Range("M80").Value = 38.82354999999995
$ws.Range("N80").Value = -2363.73685
$ws.Range("H83").Value = 647.0278
$ws.Range("I83").Value = 959.17645
$ws.Range("J83").Value = 367.73685
$ws.Range("K83").Value = 4795.882250000001
$ws.Range("L83").Value = 1838.68425
$ws.Range("M83").Value = 196.1177499999994
$ws.Range("N83").Value = -11822.68425
$ws.Range("H86").Value = 1891.9688
$ws.Range("I86").Value = 1755.625
$ws.Range("J86").Value = 2028.3125
$ws.Range("K86").Value = 1755.625
$ws.Range("L86").Value = 2028.3125
$ws.Range("M86").Value = -632.625
$ws.Range("N86").Value = -4274.3125
$ws.Range("H89").Value = 1891.9688
$ws.Range("I89").Value = 1755.625
$ws.Range("J89").Value = 2028.3125
$ws.Range("K89").Value = 8778.125
$ws.Range("L89").Value = 10141.5625
$ws.Range("M89").Value = -3162.125
$ws.Range("N89").Value = -21373.5625
$ws.Range("H107").Value = 1572.2858
$ws.Range("I107").Value = 1475
$ws.Range("J107").Value = 1815.5
$ws.Range("K107").Value = 1475
$ws.Range("L107").Value = 1815.5
$ws.Range("M107").Value = 445
$ws.Range("N107").Value = -5655.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1888.0588
$ws.Range("I31").Value = 1412.9166
$ws.Range("J31").Value = 3028.4
$ws.Range("K31").Value = 1412.9166
$ws.Range("L31").Value = 3028.4
$ws.Range("M31").Value = -1117.9166
$ws.Range("N31").Value = -3618.4
$ws.Range("H34").Value = 1888.0588
$ws.Range("I34").Value = 1412.9166
$ws.Range("J34").Value = 3028.4
$ws.Range("K34").Value = 1412.9166
$ws.Range("L34").Value = 3028.4
$ws.Range("M34").Value = -1210.9166
$ws.Range("N34").Value = -3432.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 579.36365
$ws.Range("I107").Value = 770.6429000000001
$ws.Range("J107").Value = 244.625
$ws.Range("K107").Value = 2311.9287
$ws.Range("L107").Value = 733.875
$ws.Range("M107").Value = -391.9287000000004
$ws.Range("N107").Value = -4573.875
$ws.Range("H123").Value = 4860
$ws.Range("I123").Value = 3000
$ws.Range("J123").Value = 5325
$ws.Range("K123").Value = 9000
$ws.Range("L123").Value = 15975
$ws.Range("M123").Value = -6550
$ws.Range("N123").Value = -20875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2292.3
$ws.Range("I80").Value = 2368.5715
$ws.Range("J80").Value = 2114.3333
$ws.Range("K80").Value = 2368.5715
$ws.Range("L80").Value = 2114.3333
$ws.Range("M80").Value = -1370.5715
$ws.Range("N80").Value = -4110.3333
$ws.Range("H83").Value = 2292.3
$ws.Range("I83").Value = 2368.5715
$ws.Range("J83").Value = 2114.3333
$ws.Range("K83").Value = 11842.8575
$ws.Range("L83").Value = 10571.6665
$ws.Range("M83").Value = -6850.8575
$ws.Range("N83").Value = -20555.6665
$ws.Range("H126").Value = 166668000
$ws.Range("J126").Value = 1625
$ws.Range("L126").Value = 4875
$ws.Range("N126").Value = -9815

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4670.75
$ws.Range("J61").Value = 5000
$ws.Range("L61").Value = 5000
$ws.Range("N61").Value = -5404
$ws.Range("H82").Value = 1090.2222
$ws.Range("I82").Value = 914
$ws.Range("J82").Value = 2500
$ws.Range("K82").Value = 914
$ws.Range("L82").Value = 2500
$ws.Range("M82").Value = -553
$ws.Range("N82").Value = -3222
$ws.Range("H85").Value = 1090.2222
$ws.Range("I85").Value = 914
$ws.Range("J85").Value = 2500
$ws.Range("K85").Value = 914
$ws.Range("L85").Value = 2500
$ws.Range("M85").Value = 334
$ws.Range("N85").Value = -4996
$ws.Range("H113").Value = 4670.75
$ws.Range("J113").Value = 5000
$ws.Range("L113").Value = 5000
$ws.Range("N113").Value = -9340
$ws.Range("H132").Value = 6054.263
$ws.Range("I132").Value = 6335.2
$ws.Range("J132").Value = 5000.75
$ws.Range("K132").Value = 19005.6
$ws.Range("L132").Value = 15002.25
$ws.Range("M132").Value = -16475.6
$ws.Range("N132").Value = -20062.25
$ws.Range("H136").Value = 1486.9524
$ws.Range("I136").Value = 1515.1428
$ws.Range("J136").Value = 1430.5714
$ws.Range("K136").Value = 4545.428400000001
$ws.Range("L136").Value = 4291.7142
$ws.Range("M136").Value = -1995.428400000001
$ws.Range("N136").Value = -9391.7142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1228302.8
$ws.Range("I132").Value = 2031460.1
$ws.Range("K132").Value = 6094380.300000001
$ws.Range("M132").Value = -6091850.300000001
